$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A39").Value = "Federico Fasanelli"
$ws.Range("B39").Value = "Thomas Debiasi | MAI UNA GIOIA"
$ws.Range("C39").Value = "Luca Frasca | Clitoriders"
$ws.Range("D39").Value = "Nadir Chtioui | MAI UNA GIOIA"
$ws.Range("E39").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("F39").Value = "Giacomo  Bongiovanni | Herta Vernello"
